$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DuDad-bottom-pos")

$ws.Range("B2").Value = 55.024999999999999
$ws.Range("C2").Value = -93.6
$ws.Range("D2").Value = 90

$ws.Range("B3").Value = 52.125
$ws.Range("C3").Value = -120.825

$ws.Range("B4").Value = 55.325000000000003

$ws.Range("B5").Value = 55.35
$ws.Range("C5").Value = -79.674999999999997
$ws.Range("D5").Value = 90

$ws.Range("B6").Value = 52.8
$ws.Range("C6").Value = -134.82499999999999
$ws.Range("D6").Value = 0

$ws.Range("B7").Value = 57.575000000000003

$ws.Range("B8").Value = 51.274999999999999

$ws.Range("B9").Value = 52.2

$ws.Range("B10").Value = 51.85

$ws.Range("B11").Value = 51.225000000000001

$ws.Range("B12").Value = 56.274999999999999
$ws.Range("C12").Value = -84.575000000000003
$ws.Range("D12").Value = 90

$ws.Range("B13").Value = 51.85

$ws.Range("B14").Value = 53.15

$ws.Range("B15").Value = 58.575000000000003

$ws.Range("B16").Value = 54.875

$ws.Range("D3").Select()
